$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AU2").Value = 60
$ws.Range("AV2").Value = "T"

$ws.Range("AU4").Value = 30
$ws.Range("AV4").Value = "R"
$ws.Range("AY4").Value = 90
$ws.Range("AZ4").Value = "T"

$ws.Range("AU5").Value = 60
$ws.Range("AV5").Value = "T"
$ws.Range("AY5").Value = 45
$ws.Range("AZ5").Value = "R"

$ws.Range("AU6").Value = 45
$ws.Range("AV6").Value = "R"
$ws.Range("AY6").Value = 45
$ws.Range("AZ6").Value = "T"

$ws.Range("AU7").Value = 30
$ws.Range("AV7").Value = "R"
$ws.Range("AY7").Value = 60
$ws.Range("AZ7").Value = "T"

$ws.Range("AU8").Value = 30
$ws.Range("AV8").Value = "R"
$ws.Range("AY8").Value = 30
$ws.Range("AZ8").Value = "R"

$ws.Range("AU9").Value = 60
$ws.Range("AV9").Value = "T"
$ws.Range("AY9").Value = 60
$ws.Range("AZ9").Value = "T"
$ws.Range("AX9").Value = 1

$ws.Range("AU10").Value = 45
$ws.Range("AV10").Value = "R"
$ws.Range("AY10").Value = 45
$ws.Range("AZ10").Value = "R"

$ws.Range("AU11").Value = 45
$ws.Range("AV11").Value = "T"
$ws.Range("AY11").Value = 45
$ws.Range("AZ11").Value = "T"

$ws.Range("AU12").Value = 60
$ws.Range("AV12").Value = "T"
$ws.Range("AY12").Value = 70
$ws.Range("AZ12").Value = "T"

$ws.Range("AU13").Value = 45
$ws.Range("AV13").Value = "T"
$ws.Range("AY13").Value = 60
$ws.Range("AZ13").Value = "T"

$ws.Range("AU14").Value = 45
$ws.Range("AV14").Value = "R"
$ws.Range("AY14").Value = 45
$ws.Range("AZ14").Value = "R"

$ws.Range("AU15").Value = 60
$ws.Range("AV15").Value = "T"
$ws.Range("AY15").Value = 60
$ws.Range("AZ15").Value = "T"
$ws.Range("AW15").Value = 1

$ws.Range("AU16").Value = 45
$ws.Range("AV16").Value = "T"
$ws.Range("AY16").Value = 60
$ws.Range("AZ16").Value = "T"

$ws.Range("AU18").Value = 45
$ws.Range("AV18").Value = "R"
$ws.Range("AY18").Value = 45
$ws.Range("AZ18").Value = "T"

$ws.Range("AU19").Value = 60
$ws.Range("AV19").Value = "T"

$ws.Range("AU20").Value = 20
$ws.Range("AV20").Value = "T"
$ws.Range("AY20").Value = 45
$ws.Range("AZ20").Value = "R"

$ws.Range("AU21").Value = 30
$ws.Range("AV21").Value = "R"
$ws.Range("AY21").Value = 20
$ws.Range("AZ21").Value = "R"

$ws.Range("AU22").Value = 60
$ws.Range("AV22").Value = "T"
$ws.Range("AY22").Value = 70
$ws.Range("AZ22").Value = "T"
$ws.Range("AW22").Value = 2

$ws.Range("AU24").Value = 45
$ws.Range("AV24").Value = "R"

$ws.Range("AY25").Value = 30
$ws.Range("AZ25").Value = "R"

$ws.Range("AU26").Value = 30
$ws.Range("AV26").Value = "R"
$ws.Range("AY26").Value = 30
$ws.Range("AZ26").Value = "R"
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("BC16").Select()
